$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.194999999999999
$ws.Range("D4").Value = -8.081999999999999
$ws.Range("B7").Value = 5.304
$ws.Range("D12").Value = -7.072
$ws.Range("B16").Value = 4.949
$ws.Range("D18").Value = -8.511000000000001
$ws.Range("D19").Value = -8.004999999999999
$ws.Range("D20").Value = -7.954000000000001
$ws.Range("B28").Value = 5.226000000000001
$ws.Range("B29").Value = 5.265
$ws.Range("D31").Value = -7.597
$ws.Range("B32").Value = 6.448
$ws.Range("B40").Value = 9.179
$ws.Range("D40").Value = -7.967000000000001
$ws.Range("D42").Value = -8.111000000000001
$ws.Range("D47").Value = -7.766
$ws.Range("D48").Value = -7.497999999999999
$ws.Range("B52").Value = 5.494999999999999
$ws.Range("B57").Value = 4.659999999999999
$ws.Range("D63").Value = -6.923
$ws.Range("D64").Value = -7.220000000000001
$ws.Range("B66").Value = 5.013
$ws.Range("D76").Value = -7.787999999999999
$ws.Range("D81").Value = -7.793000000000001
$ws.Range("D89").Value = -8.292
$ws.Range("D94").Value = -7.784000000000001
$ws.Range("B100").Value = 6.1
